# Update - 20150531 (2)
# The "Bamako" entries in column B (rows 51-56, the Bamako communes) were
# placeholder/incorrect values. They are corrected to match column A,
# i.e. B51:B56 now read "Commune I".."Commune VI" just like A51:A56.
# Once nothing references the shared string "Bamako" any more it drops
# out of the workbook's shared string table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B51").Value = $ws.Range("A51").Value2
$ws.Range("B52").Value = $ws.Range("A52").Value2
$ws.Range("B53").Value = $ws.Range("A53").Value2
$ws.Range("B54").Value = $ws.Range("A54").Value2
$ws.Range("B55").Value = $ws.Range("A55").Value2
$ws.Range("B56").Value = $ws.Range("A56").Value2

# Move the view so the previously-edited range is visible/selected, same
# as what Excel records after the author scrolled down and left the
# selection on the last edited cell.
$ws.Range("B56").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
